# Incorporate SOF (sofosbuvir) table updates into the drug resistance workbook.
# - Adds "producer" and "researchCode" columns (D, E) for existing drugs.
# - Adds a new row for sofosbuvir (SOF).
# - Slightly shifts the saved window x-position.
# - Sets a custom width on the new column D.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Nudge the saved workbook window position (cosmetic, matches author's re-save).
$win = $wb.Windows.Item(1)
$win.Left = 15940
$win.Top = 8540

# --- Headers for the two new columns ---
$ws.Range("D1").Value = "producer"
$ws.Range("E1").Value = "researchCode"

# --- Row 2: glecaprevir ---
$ws.Range("D2").Value = "Abbvie"
$ws.Range("E2").Value = "ABT-493"

# --- Column E for rows 3-5 (velpatasvir, pibrentasvir, voxilaprevir) ---
$ws.Range("E3").Value = "GS-5816"
$ws.Range("E4").Value = "ABT-530"
$ws.Range("E5").Value = "GS-9857"

# --- New row 6: sofosbuvir ---
$ws.Range("A6").Value = "sofosbuvir"
$ws.Range("B6").Value = "SOF"
$ws.Range("C6").Value = "NS5B RNA polymerase inhibitors"
$ws.Range("E6").Value = "GS-7977"

# --- Fill in remaining producer (column D) cells ---
$ws.Range("D3").Value = "Gilead Sciences"
$ws.Range("D4").Value = "Abbvie"
$ws.Range("D5").Value = "Gilead Sciences"
$ws.Range("D6").Value = "Gilead Sciences"

# --- Column width for the new producer column ---
$ws.Range("D1").EntireColumn.ColumnWidth = 19.29
